# tdf#90672: PPTX table cell border color is not exported.
#
# The original single-cell "Table 1" ("LibreOffice") is replaced with a
# wider two-column table ("Red" / "Green") that is repositioned/resized
# and whose cell borders + cell backgrounds are explicitly coloured, to
# exercise the (newly supported) table-cell-border-color export path.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Drop the old 1x1 "Table 1" and create a fresh table --------------
# Two scratch AddTable()/Delete() round-trips reproduce the target
# fixture's shape numbering (id=4, name="Table 3").
$old = $s.Shapes.Item(1)
$old.Delete()
$scratch = $s.Shapes.AddTable(1, 1, 0, 0, 1, 1)
$scratch.Delete()

$shape = $s.Shapes.AddTable(1, 2, 204, 210, 258, 29.2)
$shape.Name = "Table 3"

$tbl = $shape.Table
$tbl.Columns.Item(1).Width = 138
$tbl.Columns.Item(2).Width = 120
$tbl.Rows.Item(1).Height = 29.2

# --- Cell (1,1): "Red" on a light-blue border / red fill --------------
$c1 = $tbl.Cell(1, 1)
$c1.Shape.TextFrame.TextRange.Text = "Red"

foreach ($idx in 1, 2, 3) {
    # 1 = top, 2 = left, 3 = bottom
    $bd = $c1.Borders.Item($idx)
    $bd.DashStyle = 1
    $bd.Weight = 6
    $bd.ForeColor.RGB = 15773696   # 00B0F0
}
$c1.Shape.Fill.ForeColor.RGB = 255  # FF0000

# --- Cell (1,2): "Green" on a light-blue border / green fill ----------
$c2 = $tbl.Cell(1, 2)
$c2.Shape.TextFrame.TextRange.Text = "Green"

foreach ($idx in 1, 4, 3) {
    # 1 = top, 4 = right, 3 = bottom
    $bd = $c2.Borders.Item($idx)
    $bd.DashStyle = 1
    $bd.Weight = 6
    $bd.ForeColor.RGB = 15773696   # 00B0F0
}
$c2.Shape.Fill.ForeColor.RGB = 5287936  # 00B050
